# Updates cryptos list data (Price and Volume(1h) columns) to match the
# latest scrape, per commit "Updated cryptos list on Tue May  2 17:24:37 UTC 2023".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as text, preventing Excel from auto-converting
# numeric-looking strings (e.g. "22.29", "1.005") into numbers, which would
# silently drop significant trailing zeros / change the stored type.
# The NumberFormat is reset back to the workbook default ("Normal" style)
# immediately afterwards so no new cell formatting is introduced.
function Set-CellText($address, $text) {
    $range = $ws.Range($address)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-CellText "D2" "28.605.39"
Set-CellText "E2" "  +0.87%  "
Set-CellText "D3" "1.864.19"
Set-CellText "E3" "  +1.27%  "
Set-CellText "E4" "  +0.36%  "
Set-CellText "D5" "326.45"
Set-CellText "E5" "  -1.31%  "
Set-CellText "D6" "1.005"
Set-CellText "E6" "  +0.46%  "
Set-CellText "D7" "0.4621"
Set-CellText "E7" "  +0.32%  "
Set-CellText "D8" "0.3906"
Set-CellText "E8" "  +1.12%  "
Set-CellText "D9" "0.07894"
Set-CellText "E9" "  +0.29%  "
Set-CellText "D10" "0.9682"
Set-CellText "E10" "  -0.02%  "
Set-CellText "D11" "22.29"
Set-CellText "E11" "  +1.59%  "
Set-CellText "D12" "1.963.22"
Set-CellText "E12" "  +6.05%  "
Set-CellText "D13" "5.721"
Set-CellText "E13" "  -0.05%  "
Set-CellText "D14" "6.926"
Set-CellText "E14" "  -0.15%  "
Set-CellText "D15" "0.06962"
Set-CellText "E15" "  +1.12%  "
Set-CellText "D16" "88.23"
Set-CellText "E16" "  +1.40%  "
Set-CellText "D17" "1.006"
Set-CellText "E17" "  +0.52%  "
Set-CellText "E18" "  +0.91%  "
Set-CellText "E19" "  +0.16%  "
Set-CellText "D20" "1.005"
Set-CellText "E20" "  +0.38%  "
Set-CellText "D21" "28.650.68"
Set-CellText "E21" "  +0.93%  "
Set-CellText "D22" "5.310"
Set-CellText "E22" "  -0.58%  "
Set-CellText "D23" "11.06"
Set-CellText "E23" "  +0.53%  "
Set-CellText "D24" "2.125"
Set-CellText "E24" "  -1.32%  "
Set-CellText "D25" "2.157.96"
Set-CellText "E25" "  +4.35%  "
Set-CellText "D26" "153.81"
Set-CellText "E26" "  +0.09%  "
Set-CellText "D27" "19.30"
Set-CellText "E27" "  +0.37%  "
Set-CellText "D28" "5.716"
Set-CellText "E28" "  -1.26%  "
Set-CellText "D29" "1.994"
Set-CellText "E29" "  +0.39%  "
Set-CellText "D30" "119.18"
Set-CellText "E30" "  +1.97%  "
Set-CellText "E31" "  +0.25%  "
Set-CellText "D32" "0.9309"
Set-CellText "E32" "  -1.44%  "
Set-CellText "D33" "5.316"
Set-CellText "E33" "  +0.44%  "
Set-CellText "D34" "1.341"
Set-CellText "E34" "  +0.94%  "
Set-CellText "D35" "3.358"
Set-CellText "E35" "  -2.48%  "
Set-CellText "D36" "0.05830"
Set-CellText "E36" "  -3.40%  "
Set-CellText "E37" "  -1.63%  "
Set-CellText "E38" "  -0.21%  "
Set-CellText "D39" "7.886"
Set-CellText "E39" "  +3.69%  "
Set-CellText "D40" "0.5649"
Set-CellText "E40" "  +0.23%  "
Set-CellText "D41" "9.914"
Set-CellText "E41" "  -1.05%  "
Set-CellText "D42" "0.1777"
Set-CellText "E42" "  -0.43%  "
Set-CellText "D43" "0.07235"
Set-CellText "E43" "  +2.74%  "
Set-CellText "D44" "11.77"
Set-CellText "E44" "  +0.64%  "
Set-CellText "D45" "0.5311"
Set-CellText "E45" "  +0.12%  "
Set-CellText "D46" "2.166"
Set-CellText "E46" "  -4.79%  "
Set-CellText "D47" "1.140"
Set-CellText "E47" "  -7.93%  "
Set-CellText "D48" "1.845"
Set-CellText "E48" "  +0.01%  "
Set-CellText "D49" "113.21"
Set-CellText "E50" "  +0.47%  "
Set-CellText "D51" "2.342"
Set-CellText "E51" "  +0.60%  "
